{"js": "// Update the date line and every two-digit-divided-by-one-digit problem\n// in the practice sheet to the new day's values.\nconst replacements = [\n  [\"2025-10-27 Monday\", \"2025-10-28 Tuesday\"],\n  [\"14\u00f78=\", \"91\u00f78=\"],\n  [\"93\u00f73=\", \"40\u00f73=\"],\n  [\"25\u00f72=\", \"77\u00f74=\"],\n  [\"75\u00f78=\", \"15\u00f78=\"],\n  [\"89\u00f78=\", \"56\u00f79=\"],\n  [\"42\u00f77=\", \"22\u00f76=\"],\n  [\"47\u00f77=\", \"81\u00f78=\"],\n  [\"17\u00f74=\", \"93\u00f75=\"],\n  [\"38\u00f74=\", \"69\u00f73=\"],\n  [\"59\u00f78=\", \"25\u00f77=\"],\n  [\"54\u00f78=\", \"28\u00f77=\"],\n  [\"81\u00f76=\", \"22\u00f76=\"],\n  [\"13\u00f74=\", \"53\u00f79=\"],\n  [\"18\u00f77=\", \"33\u00f73=\"],\n  [\"15\u00f74=\", \"69\u00f78=\"],\n  [\"51\u00f76=\", \"57\u00f77=\"],\n  [\"98\u00f75=\", \"77\u00f77=\"],\n  [\"17\u00f73=\", \"27\u00f75=\"],\n  [\"40\u00f75=\", \"30\u00f79=\"],\n  [\"92\u00f77=\", \"30\u00f72=\"],\n  [\"24\u00f73=\", \"79\u00f77=\"],\n  [\"52\u00f72=\", \"23\u00f73=\"],\n  [\"63\u00f78=\", \"73\u00f79=\"],\n  [\"22\u00f78=\", \"91\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  if (oldText === newText) continue;\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit-divided-by-one-digit problem\n# in the practice sheet to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-10-27 Monday\", \"2025-10-28 Tuesday\"),\n  @(\"14\u00f78=\", \"91\u00f78=\"),\n  @(\"93\u00f73=\", \"40\u00f73=\"),\n  @(\"25\u00f72=\", \"77\u00f74=\"),\n  @(\"75\u00f78=\", \"15\u00f78=\"),\n  @(\"89\u00f78=\", \"56\u00f79=\"),\n  @(\"42\u00f77=\", \"22\u00f76=\"),\n  @(\"47\u00f77=\", \"81\u00f78=\"),\n  @(\"17\u00f74=\", \"93\u00f75=\"),\n  @(\"38\u00f74=\", \"69\u00f73=\"),\n  @(\"59\u00f78=\", \"25\u00f77=\"),\n  @(\"54\u00f78=\", \"28\u00f77=\"),\n  @(\"81\u00f76=\", \"22\u00f76=\"),\n  @(\"13\u00f74=\", \"53\u00f79=\"),\n  @(\"18\u00f77=\", \"33\u00f73=\"),\n  @(\"15\u00f74=\", \"69\u00f78=\"),\n  @(\"51\u00f76=\", \"57\u00f77=\"),\n  @(\"98\u00f75=\", \"77\u00f77=\"),\n  @(\"17\u00f73=\", \"27\u00f75=\"),\n  @(\"40\u00f75=\", \"30\u00f79=\"),\n  @(\"92\u00f77=\", \"30\u00f72=\"),\n  @(\"24\u00f73=\", \"79\u00f77=\"),\n  @(\"52\u00f72=\", \"23\u00f73=\"),\n  @(\"63\u00f78=\", \"73\u00f79=\"),\n  @(\"22\u00f78=\", \"91\u00f76=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    if ($oldText -eq $newText) { continue }\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n\nWrite-Output \"done\"\n"}
